# Add a "read.me" data-dictionary sheet after META, change the 8 "richness == 0"
# sentinel values in column M of META from the number 0 to the text "na", and
# move the active selection on META to M63.

$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) META!M column: replace the numeric 0 placeholders with text "na"
#    (rows 2, 11, 19, 28, 36, 45, 53, 62)
# ---------------------------------------------------------------------------
$naRows = @(2, 11, 19, 28, 36, 45, 53, 62)
foreach ($r in $naRows) {
    $meta.Range("M$r").Value = "na"
}

# ---------------------------------------------------------------------------
# 2) Insert a new worksheet "read.me" right after META, with a variable /
#    description data dictionary.
# ---------------------------------------------------------------------------
$readme = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $meta)
$readme.Name = "read.me"

$rows = @(
    @("variable", "description"),
    @("comm", "community number / code (C01 and C10 are controls)"),
    @("round", "replicate round 0, 2, 3, 4"),
    @("cyto_d00", "dilution factor to read sample on cytometer on day 00"),
    @("cyto_d01", "dilution factor to read sample on cytometer on day 01"),
    @("cyto_d02", "dilution factor to read sample on cytometer on day 02"),
    @("cyto_d07", "dilution factor to read sample on cytometer on day 07"),
    @("cyto_d14", "dilution factor to read sample on cytometer on day 14"),
    @("cyto_d21", "dilution factor to read sample on cytometer on day 21"),
    @("Sc", "S. cerevisiae present in community? (yes/no)"),
    @("Lt", "L. thermotolerans present in community? (yes/no)"),
    @("Sb", "S. bacillaris present in community? (yes/no)"),
    @("Td", "T. delbrueckii present in community? (yes/no)"),
    @("richness", "species richness of resident community (i.e., excluding B. bruxellensis and L. plantarum)")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $readme.Range("A$r").Value = $rows[$i][0]
    $readme.Range("B$r").Value = $rows[$i][1]
}

# Header row (A1:B1) gets a bold font plus the boxed border already used
# elsewhere in the workbook (reuse it via a borders-only paste so the engine
# folds onto the existing border style instead of minting a new one).
$headerRange = $readme.Range("A1:B1")
$borderSource = $meta.Range("J2")
$borderSource.Copy()
$headerRange.PasteSpecial(-4122) | Out-Null   # xlPasteBorders
$headerRange.Font.Bold = $true

# Column B is wide enough to show the full description text.
$readme.Columns.Item(2).ColumnWidth = 74.1

# View: 160% zoom, active cell B3.
$readme.Application.ActiveWindow.Zoom = 160
$readme.Range("B3").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3) Restore META as the active sheet/selection (matches the saved view in
#    the target workbook).
# ---------------------------------------------------------------------------
$meta.Activate()
$meta.Range("M63").Select() | Out-Null
